$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date column (A2:A7): "2025-12-03" -> "2025-12-05" ---
# Force text entry (avoid Excel's auto date-serial coercion), then clear the
# temporary number-format override so the cell keeps its original (default)
# style, exactly like the rest of the sheet.
$ws.Range("A2:A7").NumberFormat = "@"
$ws.Range("A2:A7").Value = "2025-12-05"
$ws.Range("A2:A7").ClearFormats()

# --- Stock identity columns (B, C) ---
$ws.Range("B2").Value = "SamsungElec"
$ws.Range("C2").Value = "005930.KS"

$ws.Range("B3").Value = "058470.KS,0P0000ASU1,98886"
$ws.Range("C3").Value = "058470.KS"

$ws.Range("B4").Value = "403870.KS,0P0001PE9K,566428"
$ws.Range("C4").Value = "403870.KS"

$ws.Range("B5").Value = "SK hynix"
$ws.Range("C5").Value = "000660.KS"

$ws.Range("B6").Value = "DB HiTek"
$ws.Range("C6").Value = "000990.KS"

$ws.Range("B7").Value = "240810.KS,0P00017YB3,330568"
$ws.Range("C7").Value = "240810.KS"

# --- Row 2 (SamsungElec) ---
$ws.Range("D2").Value = 105100
$ws.Range("E2").Value = 61.7
$ws.Range("F2").Value = 1.55
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 56
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 56.1
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 52.43913937059539
$ws.Range("O2").Value = "⚪ 중립 구간"

# --- Row 3 (058470.KS) ---
$ws.Range("D3").Value = 66000
$ws.Range("E3").Value = 67.09999999999999
$ws.Range("F3").Value = 2.33
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 56
$ws.Range("I3").Value = 56
$ws.Range("J3").Value = 66
$ws.Range("K3").Value = 56.1
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 52.43913937059539
$ws.Range("O3").Value = "⚪ 중립 구간"

# --- Row 4 (403870.KS) ---
$ws.Range("D4").Value = 31000
$ws.Range("E4").Value = 41.2
$ws.Range("F4").Value = 7.64
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 53
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 51.9
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 52.43913937059539
$ws.Range("O4").Value = "⚪ 중립 구간"

# --- Row 5 (SK hynix) ---
$ws.Range("D5").Value = 542000
$ws.Range("E5").Value = 46.2
$ws.Range("F5").Value = -0.37
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 63
$ws.Range("J5").Value = 70
$ws.Range("K5").Value = 46.9
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 52.43913937059539
$ws.Range("O5").Value = "⚪ 중립 구간"

# --- Row 6 (DB HiTek) ---
$ws.Range("D6").Value = 64400
$ws.Range("E6").Value = 41.9
$ws.Range("F6").Value = 1.26
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 46
$ws.Range("I6").Value = 50
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 44.7
$ws.Range("M6").Value = "⛔ 관망하십시오."
$ws.Range("N6").Value = 52.43913937059539
$ws.Range("O6").Value = "⚪ 중립 구간"

# --- Row 7 (240810.KS) ---
$ws.Range("D7").Value = 61900
$ws.Range("E7").Value = 32.5
$ws.Range("F7").Value = 0.98
$ws.Range("G7").Value = 20
$ws.Range("H7").Value = 66
$ws.Range("I7").Value = 56
$ws.Range("J7").Value = 56
$ws.Range("K7").Value = 44.1
$ws.Range("M7").Value = "⛔ 관망하십시오."
$ws.Range("N7").Value = 52.43913937059539
$ws.Range("O7").Value = "⚪ 중립 구간"
